# Working IEEE reader API.
# Adds s_base_mva, g_shunt_pu and b_shunt_pu columns to the "loads" sheet,
# and fixes up the v_nom_kv column (previously mislabeled/using the
# per-unit voltage instead of the nominal kV rating).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("loads")

# --- Header row ---------------------------------------------------------
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "v_nom_kv"
$ws.Range("C1").Value = "s_base_mva"
$ws.Range("D1").Value = "v_nom_pu"
$ws.Range("E1").Value = "p_nom_mw"
$ws.Range("F1").Value = "q_nom_mvar"
$ws.Range("G1").Value = "bus_idx"
$ws.Range("H1").Value = "g_shunt_pu"
$ws.Range("I1").Value = "b_shunt_pu"

# --- Data rows -----------------------------------------------------------
# name, v_nom_kv, s_base_mva, v_nom_pu, p_nom_mw, q_nom_mvar, bus_idx, g_shunt_pu, b_shunt_pu
$rows = @(
    @("Load 1", 132, 100, 1, 50, 20, 1, 0, 0),
    @("Load 2", 132, 100, 1, 50, 50, 2, 0, 0),
    @("Load 3", 132, 100, 1, 50, 20, 3, 0, 0),
    @("Load 4", 132, 100, 1, 50, 0,  4, 0, 0),
    @("Load 5", 132, 100, 1, 50, 0,  5, 0, 0)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r++
}

# --- Selection / active sheet -------------------------------------------
# Focus moves back to "busbars" (the first tab), leaving a pending
# selection of I7 on the "loads" sheet (just below the new data).
$ws.Range("I7").Select() | Out-Null

$wsBus = $wb.Worksheets.Item("busbars")
$wsBus.Activate() | Out-Null
$wsBus.Range("E8").Select() | Out-Null
